# RF014 - Administrador de Dados (Complete Test Suite): wording tweaks, v1.4 -> v1.5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text pairs (old -> new) to replace wherever they are found in the sheet.
# - Precondition: drop the "; e," -> plain "e"
# - Step 2 Expected Results: "Catalogo (Perfis) de Competencias cadastradas ... listadas"
#     -> "Perfis de Competencias cadastrados ... exibidas"
# - Step 4 Steps: drop the redundant "do novo Gerente de Desempenho" clause and trailing space
# - Step 5 Expected Results (TC1 only): drop "cadastradas"
$replacements = @{
    "Administrador esta autenticado no sistema; e, tem permissao para alterar Gerente de Desempenho" = "Administrador esta autenticado no sistema e tem permissao para alterar Gerente de Desempenho"
    "SYSTEM exibe a listagem do Catalogo (Perfis) de Competencias cadastradas com a opcao 'Alterar Gerente' dentre as varias listadas" = "SYSTEM exibe a listagem dos Perfis de Competencias cadastrados com a opcao 'Alterar Gerente' dentre as varias exibidas"
    "Administrador preenche o campo 'Login do Novo Gerente de Desempenho' do novo Gerente de Desempenho para o Perfil de Competencias " = "Administrador preenche o campo 'Login do Novo Gerente de Desempenho' para o Perfil de Competencias"
    "SYSTEM apresenta o Catalogo (Perfis) de Competencias cadastradas sem nenhuma alteracao" = "SYSTEM apresenta o Catalogo (Perfis) de Competencias sem nenhuma alteracao"
}

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
